# Backup.xlsx — Crowdin localization update
#
# Adds an English ("E" column) translation next to the existing
# Key (A) / Japanese (B) pairs in the lower lookup table (rows 26-43),
# and appends two brand-new localization keys (BackupOfficer, AllDismissItem)
# as rows 44 and 46 (row 45 stays blank, matching the sheet's existing
# blank-separator convention).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- existing rows ---------------------------------------------------
# Re-assert column B (Japanese) too: the source cells carried legacy
# <rPh>/<phoneticPr> furigana annotations (Crowdin's export drops these),
# so every B cell here is rewritten with the clean string alongside the
# new English (column E) translation.
$ws.Cells.Item(26, 2).Value = "検視官"
$ws.Cells.Item(26, 5).Value = "Coroner"

$ws.Cells.Item(27, 2).Value = "この付近に~r~死体~s~はありません。"
$ws.Cells.Item(27, 5).Value = "There is no ~r~dead bodies~s~ nearby you."

$ws.Cells.Item(28, 2).Value = "~b~{0}~s~の応援を要請しました。"
$ws.Cells.Item(28, 5).Value = "Requested ~b~{0}~s~ unit to Dispatch."

$ws.Cells.Item(29, 2).Value = "詳しい情報は~b~検視官レポート~s~を確認してください。"
$ws.Cells.Item(29, 5).Value = "You can check ~b~Coroner's Report~s~ for more information."

$ws.Cells.Item(30, 2).Value = "それではいい一日を!"
$ws.Cells.Item(30, 5).Value = "Have a nice day! Officer!"

$ws.Cells.Item(31, 2).Value = "{0}で応援を近くにテレポートさせます。"
$ws.Cells.Item(31, 5).Value = "Press {0} to teleport the backup unit nearby."

$ws.Cells.Item(33, 2).Value = "検視官メニュー"
$ws.Cells.Item(33, 5).Value = "Coroner Menu"

$ws.Cells.Item(34, 2).Value = "検視官レポート"
$ws.Cells.Item(34, 5).Value = "Coroner Report"

$ws.Cells.Item(35, 2).Value = "検視官レポート数: {0}"
$ws.Cells.Item(35, 5).Value = "Report Count: {0}"

$ws.Cells.Item(36, 2).Value = "データなし"
$ws.Cells.Item(36, 5).Value = "No Data"

$ws.Cells.Item(38, 2).Value = "名前"
$ws.Cells.Item(38, 5).Value = "Name"

$ws.Cells.Item(39, 2).Value = "性別"
$ws.Cells.Item(39, 5).Value = "Sex"

$ws.Cells.Item(40, 2).Value = "死因"
$ws.Cells.Item(40, 5).Value = "Cause of Death"

$ws.Cells.Item(41, 2).Value = "死亡日"
$ws.Cells.Item(41, 5).Value = "Died Day"

$ws.Cells.Item(43, 2).Value = "応援車両"
$ws.Cells.Item(43, 5).Value = "Backup Vehicle"

# --- brand-new rows ------------------------------------------------------
# Row 44: BackupOfficer (no English translation yet, same as source data)
$ws.Cells.Item(44, 1).Value = "BackupOfficer"
$ws.Cells.Item(44, 2).Value = "応援警官"

# Row 46: AllDismissItem
$ws.Cells.Item(46, 1).Value = "AllDismissItem"
$ws.Cells.Item(46, 2).Value = "~r~全応援を解散~s~"
$ws.Cells.Item(46, 5).Value = "~r~All Units Dismiss~s~"
